# Atualização automática dos dados - aba "Entrada"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entrada")

# Insert 6 new rows above the current row 1, pushing the existing
# rows 1-9 down to rows 7-15 (old dimension A1:F9 -> new A1:F15).
$ws.Rows("1:6").Insert()

# The "Insert" above carried the old header formatting (bold font,
# border, centered/top alignment) down with its content to row 7.
# The highlighted header look belongs on the new top row instead, so
# move that formatting from row 7 up to row 1, then clear it off row 7.
$ws.Range("A7:F7").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A7:F7").ClearFormats()

# Row 1 (new top row) - REFUGO REAL (PROCESSO)
$ws.Range("A1").Value = "REFUGO REAL (PROCESSO)"
$ws.Range("B1").Value = "R$ 236.618,39"
$ws.Range("C1").Value = "R$ 0,00"
$ws.Range("D1").Value = "R$ 236.618,39"
$ws.Range("E1").Value = "R$ 236.618,39"
$ws.Range("F1").Value = "100,00 %"

# Row 2 (new) - MATERIA PRIMA
$ws.Range("A2").Value = "MATERIA PRIMA"
$ws.Range("B2").Value = "R$ 234.121,48"
$ws.Range("C2").Value = "R$ 0,00"
$ws.Range("D2").Value = "R$ 234.121,48"
$ws.Range("E2").Value = "R$ 1,00"
$ws.Range("F2").Value = "23.412.148,00 %"

# Row 3 (new) - FRETES
$ws.Range("A3").Value = "FRETES"
$ws.Range("B3").Value = "R$ 179.648,09"
$ws.Range("C3").Value = "R$ 0,00"
$ws.Range("D3").Value = "R$ 179.648,09"
$ws.Range("E3").Value = "R$ 376.000,00"
$ws.Range("F3").Value = "47,78 %"

# Row 4 (new) - REFUGO MP+CP*
$ws.Range("A4").Value = "REFUGO MP+CP*"
$ws.Range("B4").Value = "R$ 141.770,84"
$ws.Range("C4").Value = "R$ 0,00"
$ws.Range("D4").Value = "R$ 141.770,84"
$ws.Range("E4").Value = "R$ 285.000,00"
$ws.Range("F4").Value = "49,74 %"

# Row 5 (new) - MANUTENCAO
$ws.Range("A5").Value = "MANUTENCAO"
$ws.Range("B5").Value = "R$ 121.296,45"
$ws.Range("C5").Value = "R$ 269.190,05"
$ws.Range("D5").Value = "R$ 390.486,50"
$ws.Range("E5").Value = "R$ 480.000,00"
$ws.Range("F5").Value = "81,35 %"

# Row 6 (new) - OLEOS E LUBRIFICANTES
$ws.Range("A6").Value = "OLEOS E LUBRIFICANTES"
$ws.Range("B6").Value = "R$ 80.508,54"
$ws.Range("C6").Value = "R$ 108.573,04"
$ws.Range("D6").Value = "R$ 189.081,58"
$ws.Range("E6").Value = "R$ 280.000,00"
$ws.Range("F6").Value = "67,53 %"

# Row 7 (was row 1) - CUSTO DESENVOLVIMENTO - unchanged, already shifted down

# Row 8 (was row 2) - EMBALAGENS - unchanged, already shifted down

# Row 9 (was row 3) - DESP. INDUSTRIAL - refreshed totals
$ws.Range("A9").Value = "DESP. INDUSTRIAL"
$ws.Range("B9").Value = "R$ 40.002,66"
$ws.Range("C9").Value = "R$ 100.281,50"
$ws.Range("D9").Value = "R$ 140.284,16"
$ws.Range("E9").Value = "R$ 470.000,00"
$ws.Range("F9").Value = "29,85 %"

# Row 10 (was row 4) - SERVICOS DE QUALIDADE - unchanged, already shifted down
# Row 11 (was row 5) - CUSTO COM DESENVOLVIMENTO - unchanged, already shifted down
# Row 12 (was row 6) - FERRAMENTARIA/MAN FR - unchanged, already shifted down
# Row 13 (was row 7) - MATERIAL QUALIDADE - unchanged, already shifted down
# Row 14 (was row 8) - ENERGIA ELETRICA - unchanged, already shifted down

# Row 15 (was row 9) - Total Geral - refreshed totals
$ws.Range("A15").Value = "Total Geral"
$ws.Range("B15").Value = "R$ 1.767.242,29"
$ws.Range("C15").Value = "R$ 534.600,42"
$ws.Range("D15").Value = "R$ 2.301.842,71"
$ws.Range("E15").Value = "R$ 3.835.922,27"
$ws.Range("F15").Value = "60,01 %"

Write-Output "edit applied"
